$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 1.53
$ws.Range("H3").Value = 3.75
$ws.Range("J3").Value = 2.2
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("AC3").Value = 7.5
$ws.Range("AD3").Value = 7.5
$ws.Range("AO3").Value = 8
$ws.Range("AT3").Value = 2.5
$ws.Range("M12").Value = 1.02
$ws.Range("N12").Value = 19
$ws.Range("Q12").Value = 1.53
$ws.Range("R12").Value = 2.4
$ws.Range("H15").Value = 3.15
$ws.Range("I15").Value = 2.95
$ws.Range("J15").Value = 2.87
$ws.Range("K15").Value = 2.05
$ws.Range("L15").Value = 3.5
$ws.Range("N15").Value = 6.75
$ws.Range("O15").Value = 1.33
$ws.Range("S15").Value = 1.4
$ws.Range("T15").Value = 2.52
$ws.Range("AC15").Value = 8.75
$ws.Range("AD15").Value = 6.1
$ws.Range("AE15").Value = 14.5
$ws.Range("AI15").Value = 14.5
$ws.Range("AN15").Value = 4.15
$ws.Range("AO15").Value = 12
$ws.Range("AP15").Value = 20
$ws.Range("AQ15").Value = 50
$ws.Range("AR15").Value = 80
$ws.Range("AS15").Value = 250
$ws.Range("AT15").Value = 2.5
$ws.Range("AU15").Value = 6.9
$ws.Range("AV15").Value = 60
$ws.Range("AW15").Value = 4.8
$ws.Range("AY15").Value = 24
$ws.Range("BA15").Value = 110
$ws.Range("BB15").Value = 300
$ws.Range("N16").Value = 6.75
